# Applies the "Updated symbol list" data refresh to the crypto tracker sheet.
# Columns D (Price), E (Volume 1h) and G (Hora) hold numeric-looking text that
# Excel would otherwise reinterpret as numbers/percentages, so those ranges are
# forced to Text format ("@") before the literal values are written, preserving
# exact strings such as "300.70", "0.07630" and "2,116.77%".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B/C columns: plain text (coin names, links) - safe to assign directly.
$plainTextUpdates = [ordered]@{
    'B8' = 'GateToken'
    'C8' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'B9' = 'MXToken'
    'C9' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'B11' = 'LiechtensteinCryptoassetsExchange'
    'C11' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'B12' = 'MandalaExchangeToken'
    'C12' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'B13' = 'BitrueCoin'
    'C13' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'B14' = 'BitMartToken'
    'C14' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'B15' = 'BitForexToken'
    'C15' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'B16' = 'TigerCash'
    'C16' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'B17' = 'UpBots'
    'C17' = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
    'B18' = 'LEO'
    'C18' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
}

# D/E/G columns: numeric-looking text that must stay literal text.
$forcedTextUpdates = [ordered]@{
    'D2' = '300.70'
    'E2' = '-0.03%'
    'G2' = '11'
    'D3' = '32.30'
    'E3' = '2.09%'
    'G3' = '11'
    'D4' = '4.949'
    'E4' = '-2.90%'
    'G4' = '11'
    'D5' = '0.07630'
    'E5' = '-2.16%'
    'G5' = '11'
    'D6' = '1.956'
    'E6' = '-12.63%'
    'G6' = '11'
    'D7' = '7.839'
    'E7' = '0.61%'
    'G7' = '11'
    'D8' = '3.796'
    'E8' = '-0.92%'
    'G8' = '11'
    'D9' = '0.9168'
    'E9' = '0.29%'
    'G9' = '11'
    'D10' = '0.1748'
    'E10' = '-0.17%'
    'G10' = '11'
    'D11' = '0.07786'
    'E11' = '3.22%'
    'G11' = '11'
    'D12' = '0.08539'
    'E12' = '-4.53%'
    'G12' = '11'
    'D13' = '0.03163'
    'E13' = '2.58%'
    'G13' = '11'
    'D14' = '0.09994'
    'E14' = '-0.27%'
    'G14' = '11'
    'D15' = '0.001514'
    'E15' = '-0.15%'
    'G15' = '11'
    'D16' = '0.005944'
    'E16' = '-0.88%'
    'G16' = '11'
    'D17' = '0.007498'
    'E17' = '2,116.77%'
    'G17' = '11'
    'D18' = '3.461'
    'E18' = '-0.12%'
    'G18' = '11'
    'D19' = '2.153'
    'E19' = '-4.44%'
    'G19' = '11'
    'E20' = '1.82%'
    'G20' = '11'
    'D21' = '0.1300'
    'E21' = '-2.78%'
    'G21' = '11'
    'D22' = '4.265'
    'E22' = '5.00%'
    'G22' = '11'
    'D23' = '0.1993'
    'E23' = '9.68%'
    'G23' = '11'
    'D24' = '0.04514'
    'E24' = '-1.58%'
    'G24' = '11'
    'D25' = '0.001221'
    'E25' = '-2.37%'
    'G25' = '11'
    'D26' = '0.004386'
    'E26' = '-1.79%'
    'G26' = '11'
    'D27' = '0.0001251'
    'E27' = '0.11%'
    'G27' = '11'
    'G28' = '11'
    'G29' = '11'
    'G30' = '11'
    'G31' = '11'
    'G32' = '11'
    'G33' = '11'
    'G34' = '11'
    'G35' = '11'
    'G36' = '11'
    'G37' = '11'
    'G38' = '11'
    'D39' = '0.01696'
    'E39' = '-4.04%'
    'G39' = '11'
    'D40' = '0.04681'
    'E40' = '-1.73%'
    'G40' = '11'
    'D41' = '0.007466'
    'E41' = '-2.89%'
    'G41' = '11'
    'D42' = '0.1348'
    'E42' = '-0.68%'
    'G42' = '11'
    'D43' = '0.002332'
    'E43' = '6.51%'
    'G43' = '11'
    'D44' = '0.01049'
    'E44' = '2.78%'
    'G44' = '11'
    'D45' = '0.00006262'
    'E45' = '1.37%'
    'G45' = '11'
    'E46' = '0.12%'
    'G46' = '11'
    'D47' = '0.8205'
    'E47' = '10.52%'
    'G47' = '11'
    'G48' = '11'
    'E49' = '0.12%'
    'G49' = '11'
    'E50' = '0.12%'
    'G50' = '11'
    'G51' = '11'
}

foreach ($ref in $plainTextUpdates.Keys) {
    $ws.Range($ref).Value = $plainTextUpdates[$ref]
}

foreach ($ref in $forcedTextUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $forcedTextUpdates[$ref]
}
